# Update SwaadSutra_Daily_2026-01-21.xlsx
# Order #27 (Renu) has moved from COOKING to READY status.

$wb = $excel.ActiveWorkbook

# --- Daily Orders sheet: update order status ---
$ordersSheet = $wb.Worksheets.Item("Daily Orders")
$ordersSheet.Range("H3").Value = "READY"

# --- Summary sheet: adjust Cooking/Ready counts to reflect the status change ---
$summarySheet = $wb.Worksheets.Item("Summary")
$summarySheet.Range("C2").Value = 0
$summarySheet.Range("D2").Value = 1
